$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team_PER_2011: fixed PER bug - re-map team codes (col B) to the correct rows
$ws.Range("B2").Value = "POR"
$ws.Range("B3").Value = "NJN"
$ws.Range("B4").Value = "CLE"
$ws.Range("B5").Value = "DAL"
$ws.Range("B6").Value = "ATL"
$ws.Range("B7").Value = "OKC"
$ws.Range("B8").Value = "CHA"
$ws.Range("B9").Value = "WAS"
$ws.Range("B10").Value = "MIL"
$ws.Range("B11").Value = "LAC"
$ws.Range("B12").Value = "SAS"
$ws.Range("B13").Value = "DET"
$ws.Range("B14").Value = "ORL"
$ws.Range("B15").Value = "UTA"
$ws.Range("B16").Value = "MEM"
$ws.Range("B17").Value = "HOU"
$ws.Range("B18").Value = "DEN"
$ws.Range("B19").Value = "LAL"
$ws.Range("B20").Value = "GSW"
$ws.Range("B21").Value = "IND"
$ws.Range("B22").Value = "CHI"
$ws.Range("B23").Value = "PHI"
$ws.Range("B24").Value = "BOS"
$ws.Range("B25").Value = "TOR"
$ws.Range("B26").Value = "MIA"
$ws.Range("B27").Value = "SAC"
$ws.Range("B28").Value = "PHO"
$ws.Range("B29").Value = "NOH"
$ws.Range("B30").Value = "NYK"
$ws.Range("B31").Value = "MIN"

# Corrected PER values in column C
$ws.Range("C2").Value = 11.725
$ws.Range("C3").Value = 13.08333333333333
$ws.Range("C4").Value = 12.05
$ws.Range("C5").Value = 14.23076923076923
$ws.Range("C6").Value = 13.34166666666667
$ws.Range("C7").Value = 12.28333333333333
$ws.Range("C8").Value = 11.66666666666667
$ws.Range("C9").Value = 11.45
$ws.Range("C10").Value = 12.42666666666667
$ws.Range("C11").Value = 13.6
$ws.Range("C12").Value = 12.83333333333334
$ws.Range("C13").Value = 14.06153846153846
$ws.Range("C14").Value = 13.275
$ws.Range("C15").Value = 12.19285714285714
$ws.Range("C16").Value = 15.3
$ws.Range("C17").Value = 16.63636363636364
$ws.Range("C18").Value = 14.85555555555556
$ws.Range("C19").Value = 14.03076923076923
$ws.Range("C20").Value = 13.53333333333333
$ws.Range("C21").Value = 12.66666666666666
$ws.Range("C22").Value = 12.9
$ws.Range("C23").Value = 11.27333333333333
$ws.Range("C24").Value = 13.63636363636364
$ws.Range("C25").Value = 11.72857142857143
$ws.Range("C26").Value = 11.9
$ws.Range("C27").Value = 9.207692307692305
$ws.Range("C28").Value = 12.34
$ws.Range("C29").Value = 11.23333333333333
$ws.Range("C30").Value = 11.8
$ws.Range("C31").Value = 12.4
